$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-26 19:18:23'
$ws.Range('E3').Value = '2026-02-26 19:18:26'
$ws.Range('O3').Value = '3.1 °C'
$ws.Range('E4').Value = '2026-02-26 19:18:28'
$ws.Range('E5').Value = '2026-02-26 19:18:31'
$ws.Range('O5').Value = '5.2 °C'
$ws.Range('E6').Value = '2026-02-26 19:18:33'
$ws.Range('E7').Value = '2026-02-26 19:18:36'
$ws.Range('E8').Value = '2026-02-26 19:18:38'
$ws.Range('E9').Value = '2026-02-26 19:18:41'
$ws.Range('O9').Value = '12.6 °C'
$ws.Range('E10').Value = '2026-02-26 19:18:44'
$ws.Range('O10').Value = '9.8 °C'
$ws.Range('E11').Value = '2026-02-26 19:18:46'
$ws.Range('E12').Value = '2026-02-26 19:18:48'
$ws.Range('H12').Value = "'91%"
$ws.Range('E13').Value = '2026-02-26 19:18:51'
$ws.Range('J13').Value = '1028.1 hPa'
$ws.Range('L13').Value = '29.9 km/h - 115º 18:45 TU'
$ws.Range('O13').Value = '7.3 °C'
$ws.Range('E14').Value = '2026-02-26 19:18:53'
$ws.Range('E15').Value = '2026-02-26 19:18:56'
$ws.Range('O15').Value = '12.0 °C'
$ws.Range('E16').Value = '2026-02-26 19:18:58'
$ws.Range('E17').Value = '2026-02-26 19:19:00'
$ws.Range('E18').Value = '2026-02-26 19:19:03'
$ws.Range('E19').Value = '2026-02-26 19:19:05'
$ws.Range('H19').Value = "'45%"
$ws.Range('K19').Value = '15.9 MJ/m2'
$ws.Range('E20').Value = '2026-02-26 19:19:07'
$ws.Range('K20').Value = '17.1 MJ/m2'
$ws.Range('O20').Value = '2.8 °C'
$ws.Range('E21').Value = '2026-02-26 19:19:10'
$ws.Range('O21').Value = '10.1 °C'
$ws.Range('E22').Value = '2026-02-26 19:19:12'
$ws.Range('E23').Value = '2026-02-26 19:19:15'
$ws.Range('E24').Value = '2026-02-26 19:19:17'
$ws.Range('E25').Value = '2026-02-26 19:19:20'
$ws.Range('E26').Value = '2026-02-26 19:19:22'
$ws.Range('J26').Value = '1024.2 hPa'
$ws.Range('O26').Value = '11.1 °C'
$ws.Range('E27').Value = '2026-02-26 19:19:25'
$ws.Range('E28').Value = '2026-02-26 19:19:27'
$ws.Range('O28').Value = '11.1 °C'
$ws.Range('E29').Value = '2026-02-26 19:19:30'
$ws.Range('O29').Value = '11.9 °C'
$ws.Range('E30').Value = '2026-02-26 19:19:32'
$ws.Range('H30').Value = "'85%"
$ws.Range('O30').Value = '12.4 °C'
$ws.Range('E31').Value = '2026-02-26 19:19:34'
$ws.Range('J31').Value = '1026.7 hPa'
$ws.Range('E32').Value = '2026-02-26 19:19:37'
$ws.Range('H32').Value = "'63%"
$ws.Range('O32').Value = '8.4 °C'
$ws.Range('E33').Value = '2026-02-26 19:19:40'
$ws.Range('J33').Value = '1026.6 hPa'
$ws.Range('O33').Value = '8.8 °C'
$ws.Range('E34').Value = '2026-02-26 19:19:42'
$ws.Range('O34').Value = '5.0 °C'
$ws.Range('E35').Value = '2026-02-26 19:19:45'
$ws.Range('J35').Value = '1025.3 hPa'
$ws.Range('O35').Value = '12.4 °C'
$ws.Range('E36').Value = '2026-02-26 19:19:47'
$ws.Range('O36').Value = '12.8 °C'
$ws.Range('E37').Value = '2026-02-26 19:19:49'
$ws.Range('H37').Value = "'72%"
$ws.Range('E38').Value = '2026-02-26 19:19:52'
$ws.Range('H38').Value = "'79%"
$ws.Range('E39').Value = '2026-02-26 19:19:54'
$ws.Range('H39').Value = "'42%"
$ws.Range('E40').Value = '2026-02-26 19:19:57'
$ws.Range('E41').Value = '2026-02-26 19:19:59'
$ws.Range('O41').Value = '11.3 °C'
$ws.Range('E42').Value = '2026-02-26 19:20:02'
$ws.Range('E43').Value = '2026-02-26 19:20:04'
$ws.Range('E44').Value = '2026-02-26 19:20:06'
$ws.Range('H44').Value = "'52%"
$ws.Range('E45').Value = '2026-02-26 19:20:09'
$ws.Range('E46').Value = '2026-02-26 19:20:11'
$ws.Range('H46').Value = "'81%"
$ws.Range('J46').Value = '1027.0 hPa'
